# "Optimized csv file mode." - update the BackTester test-case rows on Sheet1:
#  - narrow/shift the From-To date ranges used for the CSV-driven backtests
#  - switch both test rows to the MACD strategy (was EarlyMACD)
#  - switch row 2's exchange/symbol to ByBit/BTCUSDT (was Binance/ETHUSDT)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 2 (Test #1)
$ws1.Range("D2").Value = 44197
$ws1.Range("E2").Value = 44555
$ws1.Range("L2").Value = "MACD"

# Row 3 (Test #2)
$ws1.Range("B3").Value = "ByBit"
$ws1.Range("C3").Value = "BTCUSDT"
$ws1.Range("D3").Value = 44197
$ws1.Range("E3").Value = 44555
$ws1.Range("L3").Value = "MACD"

# Leave the sheet's active selection on L2, matching the saved workbook state
$ws1.Activate()
$ws1.Range("L2").Select()
